$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map of 1-indexed table row number -> array of new cell values (left to right).
# Only the rows that contain data (1, 5, 9, 13, 17) need updates; the text
# is written cell-by-cell so each w:tc keeps its existing run formatting
# (font / size) untouched, matching the diff (only <w:t> contents change).

$rowUpdates = @{
    1  = @("38÷2=", "24÷6=", "94÷9=", "67÷5=", "33÷6=")
    5  = @("17÷2=", "25÷3=", "57÷8=", "94÷9=", "66÷5=")
    9  = @("86÷8=", "12÷6=", "64÷4=", "38÷3=", "41÷8=")
    13 = @("21÷3=", "16÷3=", "44÷3=", "21÷9=", "42÷4=")
    17 = @("56÷6=", "27÷5=", "95÷8=", "28÷6=", "17÷9=")
}

foreach ($rowIndex in $rowUpdates.Keys) {
    $values = $rowUpdates[$rowIndex]
    $row = $t.Rows.Item($rowIndex)
    for ($i = 1; $i -le $values.Length; $i++) {
        $cell = $row.Cells.Item($i)
        $cell.Range.Text = $values[$i - 1]
    }
}

Write-Output "done"
